# Apply the crypto price/volume updates captured in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# (these prices are stored as literal text in the workbook, e.g. "303.59").
# Temporarily force Text format so the value lands as a string, then clear
# the formatting again so no stray number-format style is left behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$ws.Range("D2").Value = "23.257.37"
$ws.Range("E2").Value = "  +1.06%  "

$ws.Range("D3").Value = "1.604.37"
$ws.Range("E3").Value = "  +0.13%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("E5").Value = "  -0.08%  "

Set-TextValue $ws.Range("D6") "303.59"
$ws.Range("E6").Value = "  +0.80%  "

Set-TextValue $ws.Range("D7") "0.3770"
$ws.Range("E7").Value = "  -0.27%  "

$ws.Range("E8").Value = "  +4.70%  "

Set-TextValue $ws.Range("D9") "0.3640"
$ws.Range("E9").Value = "  +0.33%  "

Set-TextValue $ws.Range("D10") "1.279"
$ws.Range("E10").Value = "  +1.48%  "

Set-TextValue $ws.Range("D11") "1.001"
$ws.Range("E11").Value = "  -0.07%  "

Set-TextValue $ws.Range("D12") "0.08130"
$ws.Range("E12").Value = "  +0.02%  "

$ws.Range("E13").Value = "  +0.04%  "

Set-TextValue $ws.Range("D14") "6.607"
$ws.Range("E14").Value = "  +0.20%  "

Set-TextValue $ws.Range("D15") "7.434"
$ws.Range("E15").Value = "  +0.51%  "

Set-TextValue $ws.Range("D16") "0.00001250"
$ws.Range("E16").Value = "  +0.54%  "

$ws.Range("D17").Value = "1.603.01"
$ws.Range("E17").Value = "  +0.28%  "

Set-TextValue $ws.Range("D18") "94.06"
$ws.Range("E18").Value = "  +2.07%  "

Set-TextValue $ws.Range("D19") "0.06947"
$ws.Range("E19").Value = "  +1.02%  "

Set-TextValue $ws.Range("D20") "18.20"
$ws.Range("E20").Value = "  -0.25%  "

Set-TextValue $ws.Range("D21") "6.535"
$ws.Range("E21").Value = "  -0.45%  "

$ws.Range("E22").Value = "  -0.10%  "

$ws.Range("E23").Value = "  -1.48%  "

$ws.Range("D24").Value = "23.241.95"

Set-TextValue $ws.Range("D25") "3.059"
$ws.Range("E25").Value = "  +9.32%  "

$ws.Range("E26").Value = "  +0.69%  "

Set-TextValue $ws.Range("D27") "21.25"
$ws.Range("E27").Value = "  +0.74%  "

Set-TextValue $ws.Range("D28") "149.88"
$ws.Range("E28").Value = "  -0.31%  "

Set-TextValue $ws.Range("D29") "5.260"
$ws.Range("E29").Value = "  +0.14%  "

Set-TextValue $ws.Range("D30") "134.83"
$ws.Range("E30").Value = "  +0.86%  "

Set-TextValue $ws.Range("D31") "2.397"
$ws.Range("E31").Value = "  +3.51%  "

Set-TextValue $ws.Range("D32") "6.745"
$ws.Range("E32").Value = "  -0.93%  "

$ws.Range("D33").Value = "1.782.54"

Set-TextValue $ws.Range("D34") "0.9663"
$ws.Range("E34").Value = "  +0.44%  "

Set-TextValue $ws.Range("D35") "0.07506"
$ws.Range("E35").Value = "  -1.68%  "

$ws.Range("E36").Value = "  +1.93%  "

Set-TextValue $ws.Range("D37") "10.34"
$ws.Range("E37").Value = "  -0.50%  "

$ws.Range("E38").Value = "  +0.10%  "

$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D39") "0.08819"
$ws.Range("E39").Value = "  -0.44%  "

$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D40") "6.124"
$ws.Range("E40").Value = "  -2.76%  "

Set-TextValue $ws.Range("D41") "1.395"
$ws.Range("E41").Value = "  +2.27%  "

Set-TextValue $ws.Range("D42") "0.7128"
$ws.Range("E42").Value = "  +0.99%  "

Set-TextValue $ws.Range("D43") "12.52"
$ws.Range("E43").Value = "  +0.13%  "

Set-TextValue $ws.Range("D44") "15.58"
$ws.Range("E44").Value = "  +2.67%  "

Set-TextValue $ws.Range("D45") "0.6564"
$ws.Range("E45").Value = "  -0.86%  "

Set-TextValue $ws.Range("D46") "2.321"
$ws.Range("E46").Value = "  +0.20%  "

Set-TextValue $ws.Range("D47") "0.9994"
$ws.Range("E47").Value = "  -0.10%  "

Set-TextValue $ws.Range("D48") "4.016"
$ws.Range("E48").Value = "  +0.59%  "

Set-TextValue $ws.Range("D49") "132.65"
$ws.Range("E49").Value = "  -0.05%  "

Set-TextValue $ws.Range("D50") "0.07954"
$ws.Range("E50").Value = "  +0.62%  "

Set-TextValue $ws.Range("D51") "1.208"
$ws.Range("E51").Value = "  -1.01%  "
